$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.301.44"
$ws.Range("E2").Value = "  +0.65%  "

# Row 3
$ws.Range("D3").Value = "1.679.85"
$ws.Range("E3").Value = "  +0.75%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.40"
$ws.Range("E5").Value = "  +0.76%  "

# Row 6
$ws.Range("E6").Value = "  +2.61%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.007"
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2700"
$ws.Range("E8").Value = "  +2.56%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06479"
$ws.Range("E9").Value = "  +1.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.98"
$ws.Range("E10").Value = "  +1.73%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07523"
$ws.Range("E11").Value = "  +1.44%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.527"
$ws.Range("E12").Value = "  +0.29%  "

# Row 13
$ws.Range("D13").Value = "1.670.04"
$ws.Range("E13").Value = "  -0.14%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5804"
$ws.Range("E14").Value = "  -0.06%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008524"
$ws.Range("E15").Value = "  -0.54%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.69"
$ws.Range("E16").Value = "  +0.77%  "

# Row 17
$ws.Range("D17").Value = "26.334.79"
$ws.Range("E17").Value = "  +0.53%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.919"
$ws.Range("E18").Value = "  -0.16%  "

# Row 19
$ws.Range("E19").Value = "  +0.17%  "

# Row 20
$ws.Range("E20").Value = "  +0.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.93"
$ws.Range("E21").Value = "  +0.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.204"
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.63"
$ws.Range("E24").Value = "  +0.19%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.792"
$ws.Range("E25").Value = "  +2.15%  "

# Row 26
$ws.Range("E26").Value = "  +4.17%  "

# Row 27
$ws.Range("E27").Value = "  +1.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06463"
$ws.Range("E28").Value = "  +2.75%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.356"
$ws.Range("E29").Value = "  +4.68%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.330"
$ws.Range("E30").Value = "  +0.91%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.602"
$ws.Range("E31").Value = "  +2.13%  "

# Row 32
$ws.Range("E32").Value = "  +2.48%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.660"
$ws.Range("E33").Value = "  +1.27%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.032"
$ws.Range("E34").Value = "  +1.76%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6254"
$ws.Range("E35").Value = "  +2.96%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.403"
$ws.Range("E36").Value = "  +1.65%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.727"
$ws.Range("E37").Value = "  +3.06%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.455"
$ws.Range("E38").Value = "  +4.63%  "

# Row 39
$ws.Range("D39").Value = "1.110.57"
$ws.Range("E39").Value = "  +2.88%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01622"
$ws.Range("E40").Value = "  +0.89%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8786"
$ws.Range("E41").Value = "  +1.65%  "

# Row 42
$ws.Range("E42").Value = "  +0.48%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.70"
$ws.Range("E43").Value = "  -0.39%  "

# Row 44
$ws.Range("D44").Value = "1.832.50"
$ws.Range("E44").Value = "  +0.96%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000113"
$ws.Range("E45").Value = "  +1.56%  "

# Row 46
$ws.Range("E46").Value = "  +1.66%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.185"
$ws.Range("E47").Value = "  +1.31%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.009"
$ws.Range("E48").Value = "  +0.03%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05270"
$ws.Range("E49").Value = "  +1.33%  "

# Row 50
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4292"
$ws.Range("E50").Value = "  +0.00%  "

# Row 51
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.076"
$ws.Range("E51").Value = "  +2.94%  "
